{"js": "// Replace each \"two-digit \u00f7 one-digit\" answer cell with its updated value.\n// Each (oldText, newText) pair occurs exactly once in the document, so a\n// scoped search + full-range replace is unambiguous for every entry.\nconst replacements = [\n  [\"16\u00f76=2, 4\", \"23\u00f79=2, 5\"],\n  [\"57\u00f79=6, 3\", \"71\u00f79=7, 8\"],\n  [\"64\u00f77=9, 1\", \"73\u00f76=12, 1\"],\n  [\"79\u00f76=13, 1\", \"72\u00f73=24, 0\"],\n  [\"29\u00f78=3, 5\", \"40\u00f72=20, 0\"],\n  [\"80\u00f76=13, 2\", \"90\u00f74=22, 2\"],\n  [\"33\u00f74=8, 1\", \"68\u00f73=22, 2\"],\n  [\"57\u00f72=28, 1\", \"39\u00f79=4, 3\"],\n  [\"68\u00f78=8, 4\", \"46\u00f72=23, 0\"],\n  [\"32\u00f77=4, 4\", \"52\u00f77=7, 3\"],\n  [\"65\u00f76=10, 5\", \"57\u00f76=9, 3\"],\n  [\"14\u00f73=4, 2\", \"97\u00f74=24, 1\"],\n  [\"54\u00f79=6, 0\", \"43\u00f79=4, 7\"],\n  [\"74\u00f76=12, 2\", \"84\u00f74=21, 0\"],\n  [\"61\u00f74=15, 1\", \"10\u00f75=2, 0\"],\n  [\"41\u00f78=5, 1\", \"84\u00f74=21, 0\"],\n  [\"79\u00f74=19, 3\", \"60\u00f79=6, 6\"],\n  [\"74\u00f77=10, 4\", \"91\u00f72=45, 1\"],\n  [\"40\u00f75=8, 0\", \"44\u00f72=22, 0\"],\n  [\"77\u00f79=8, 5\", \"10\u00f73=3, 1\"],\n  [\"86\u00f76=14, 2\", \"99\u00f75=19, 4\"],\n  [\"37\u00f77=5, 2\", \"50\u00f75=10, 0\"],\n  [\"83\u00f75=16, 3\", \"98\u00f74=24, 2\"],\n  [\"39\u00f75=7, 4\", \"90\u00f74=22, 2\"],\n  [\"83\u00f79=9, 2\", \"93\u00f74=23, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each (old, new) pair occurs exactly once in the document, so a single\n# Find/Replace pass (wdReplaceOne) against the whole-document range is\n# unambiguous for every entry and leaves all other text untouched.\n$replacements = @(\n    @(\"16\u00f76=2, 4\", \"23\u00f79=2, 5\"),\n    @(\"57\u00f79=6, 3\", \"71\u00f79=7, 8\"),\n    @(\"64\u00f77=9, 1\", \"73\u00f76=12, 1\"),\n    @(\"79\u00f76=13, 1\", \"72\u00f73=24, 0\"),\n    @(\"29\u00f78=3, 5\", \"40\u00f72=20, 0\"),\n    @(\"80\u00f76=13, 2\", \"90\u00f74=22, 2\"),\n    @(\"33\u00f74=8, 1\", \"68\u00f73=22, 2\"),\n    @(\"57\u00f72=28, 1\", \"39\u00f79=4, 3\"),\n    @(\"68\u00f78=8, 4\", \"46\u00f72=23, 0\"),\n    @(\"32\u00f77=4, 4\", \"52\u00f77=7, 3\"),\n    @(\"65\u00f76=10, 5\", \"57\u00f76=9, 3\"),\n    @(\"14\u00f73=4, 2\", \"97\u00f74=24, 1\"),\n    @(\"54\u00f79=6, 0\", \"43\u00f79=4, 7\"),\n    @(\"74\u00f76=12, 2\", \"84\u00f74=21, 0\"),\n    @(\"61\u00f74=15, 1\", \"10\u00f75=2, 0\"),\n    @(\"41\u00f78=5, 1\", \"84\u00f74=21, 0\"),\n    @(\"79\u00f74=19, 3\", \"60\u00f79=6, 6\"),\n    @(\"74\u00f77=10, 4\", \"91\u00f72=45, 1\"),\n    @(\"40\u00f75=8, 0\", \"44\u00f72=22, 0\"),\n    @(\"77\u00f79=8, 5\", \"10\u00f73=3, 1\"),\n    @(\"86\u00f76=14, 2\", \"99\u00f75=19, 4\"),\n    @(\"37\u00f77=5, 2\", \"50\u00f75=10, 0\"),\n    @(\"83\u00f75=16, 3\", \"98\u00f74=24, 2\"),\n    @(\"39\u00f75=7, 4\", \"90\u00f74=22, 2\"),\n    @(\"83\u00f79=9, 2\", \"93\u00f74=23, 1\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $oldText, $false, $false, $false, $false, $false,\n        $true, 1, $false, $newText, 1\n    )\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
